$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for the Zapallo / Paine series.
# It belongs right before the current row 461, so insert a fresh row
# there; Excel shifts every row from 461..563 down to 462..564
# automatically (carrying along values, styles, etc.), which matches
# the diff exactly.
$ws.Rows.Item(461).Insert()

# Populate the newly inserted row 461 with the new record's data.
$ws.Range("A461").Value = 4
$ws.Range("B461").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C461").Value = "Los Lagos"
$ws.Range("D461").Value = 45244
$ws.Range("D461").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E461").Value = 10
$ws.Range("F461").Value = 100112045
$ws.Range("G461").Value = "Zapallo"
$ws.Range("H461").Value = "Paine"
$ws.Range("I461").Value = "1a (guarda)"
$ws.Range("J461").Value = 1200
$ws.Range("K461").Value = 1400
$ws.Range("L461").Value = 1400
$ws.Range("M461").Value = 1400
$ws.Range("N461").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O461").Value = "Región de O'Higgins"
$ws.Range("P461").Value = 1400
$ws.Range("Q461").Value = 1
$ws.Range("R461").Value = "Hortaliza"
